$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested count) values for two rows.
# These figures live on both the "展览" sheet and the aggregated "全部类型" sheet.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 544
    $ws.Range("F7").Value = 779
}
